$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: Cruise -> Month
$ws.Range("A1").Value = "Month"

# Cruise code -> Month name
$ws.Range("A2").Value = "March"
$ws.Range("A3").Value = "March"
$ws.Range("A4").Value = "March"
$ws.Range("A5").Value = "March"
$ws.Range("A6").Value = "March"

$ws.Range("A7").Value = "October"
$ws.Range("A8").Value = "October"
$ws.Range("A9").Value = "October"
$ws.Range("A10").Value = "October"
$ws.Range("A11").Value = "October"
$ws.Range("A12").Value = "October"

# Tiny last-digit floating point refresh on recomputed summary values
$ws.Range("E2").Value = 51.88829818482504
$ws.Range("F2").Value = 63.33526841717942

$ws.Range("E7").Value = 93.61407208076504
$ws.Range("F7").Value = 13.04166617688351

$ws.Range("F9").Value = 0.8088966589264726

$ws.Range("D10").Value = 777.9622122079993
$ws.Range("E10").Value = 30.04640636922772

$ws.Range("C11").Value = 6663.286950780684

$ws.Range("E12").Value = 35.38082624789594
